$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3.8
$ws.Range("K2").Value = 2.38
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.25
$ws.Range("S2").Value = 1.3
$ws.Range("T2").Value = 3.4
$ws.Range("AH2").Value = 9
$ws.Range("AT2").Value = 3.4
$ws.Range("AX2").Value = 9
